$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Valor Mora" total (E11): 1443304 -> 1478416
$ws.Cells.Item(11, 5).Value = 1478416

# 2. Update "Cant. Periodos" count (F13): 37 -> 38
$ws.Cells.Item(13, 6).Value = 38

# 3. Insert a new data row (row 54) for period "2509", right before the
#    footer/signature rows, copying the formatting of the previous last
#    data row (row 53, period "2508").
$ws.Rows.Item(54).Insert()
$ws.Range("B53:J53").Copy($ws.Range("B54:J54"))

# The new row keeps the same worker/period data pattern as the row above,
# only the period label changes to the new period "2509".
$ws.Cells.Item(54, 5).Value = "2509"

$wb.Saved = $false
